$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: remove the red "PETUNJUK" instruction call-outs in B1/C1,
#     replacing them with empty, right-aligned / wrap-text cells that keep
#     the same shaded fill used elsewhere on the sheet (copied from A1). ---
$ws.Range("B1").ClearContents()
$ws.Range("C1").ClearContents()

$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("B1:C1").HorizontalAlignment = -4152   # xlRight
$ws.Range("B1:C1").VerticalAlignment = -4107     # xlBottom (Excel default -> no explicit attr)
$ws.Range("B1:C1").WrapText = $true

# --- Data rows: turn the "[n]" bracket placeholder text into real numbers,
#     and give the "X" / number cells a cleaner alignment (center / right,
#     no explicit vertical alignment). ---
$ws.Range("A3:A4").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A3:A4").VerticalAlignment = -4107     # xlBottom (default)

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 3.1

$ws.Range("B3:C4").HorizontalAlignment = -4152   # xlRight
$ws.Range("B3:C4").VerticalAlignment = -4107     # xlBottom (default)

# --- New row 5, continuing the same "X" / code pattern. ---
$ws.Range("A5").Value = "X"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 4.2
$ws.Range("B5:C5").HorizontalAlignment = -4152   # xlRight
$ws.Range("B5:C5").VerticalAlignment = -4107     # xlBottom (default)

# --- Columns B/C default style (used by any future cell in the column)
#     now matches the new right-aligned look. ---
$ws.Columns("B:C").HorizontalAlignment = -4152   # xlRight
$ws.Columns("B:C").VerticalAlignment = -4107     # xlBottom (default)

# --- Selection moved on save, as recorded in the workbook. ---
$ws.Range("C10").Select()
